$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.Value = "'42.899.58"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$r = $ws.Range("D3")
$r.Value = "'2.535.59"
$r.Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "
$r = $ws.Range("D4")
$r.Value = "'1.00"
$r.Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$r = $ws.Range("D5")
$r.Value = "'318.67"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.70%  "
$r = $ws.Range("D6")
$r.Value = "'97.24"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  -0.05%  "
$r = $ws.Range("D9")
$r.Value = "'0.538"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "
$r = $ws.Range("D10")
$r.Value = "'35.84"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -1.62%  "
$r = $ws.Range("D11")
$r.Value = "'0.0817"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "
$r = $ws.Range("D12")
$r.Value = "'7.55"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  -3.64%  "
$r = $ws.Range("D14")
$r.Value = "'2.923.89"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -0.22%  "
$r = $ws.Range("D15")
$r.Value = "'2.630.59"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +2.24%  "
$r = $ws.Range("D16")
$r.Value = "'15.11"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -3.70%  "
$r = $ws.Range("D17")
$r.Value = "'0.853"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -1.53%  "
$r = $ws.Range("D18")
$r.Value = "'42.926.50"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +0.66%  "
$r = $ws.Range("D19")
$r.Value = "'6.92"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +3.40%  "
$r = $ws.Range("D20")
$r.Value = "'12.70"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -2.96%  "
$r = $ws.Range("D21")
$r.Value = "'0.0₃0968"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.43%  "
$r = $ws.Range("D22")
$r.Value = "'69.73"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "
$r = $ws.Range("D23")
$r.Value = "'253.28"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "
$r = $ws.Range("D24")
$r.Value = "'2.99"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +1.06%  "
$r = $ws.Range("D26")
$r.Value = "'26.47"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("E27").Value = "  +0.79%  "
$r = $ws.Range("D28")
$r.Value = "'2.42"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +2.38%  "
$r = $ws.Range("D29")
$r.Value = "'41.35"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +3.78%  "
$r = $ws.Range("D30")
$r.Value = "'10.51"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +4.28%  "
$r = $ws.Range("D31")
$r.Value = "'5.91"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "
$r = $ws.Range("D32")
$r.Value = "'157.71"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$r = $ws.Range("D34")
$r.Value = "'3.36"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$r = $ws.Range("D35")
$r.Value = "'19.27"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E36").Value = "  +3.09%  "
$r = $ws.Range("D37")
$r.Value = "'0.0792"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +11.76%  "
$r = $ws.Range("D41")
$r.Value = "'21.82"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -12.52%  "
$r = $ws.Range("D42")
$r.Value = "'0.0306"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +0.66%  "
$r = $ws.Range("D43")
$r.Value = "'3.81"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  +0.28%  "
$r = $ws.Range("D45")
$r.Value = "'3.28"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -3.53%  "
$r = $ws.Range("D46")
$r.Value = "'2.011.11"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$r = $ws.Range("D47")
$r.Value = "'9.17"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "
$r = $ws.Range("D48")
$r.Value = "'84.34"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -1.25%  "
$r = $ws.Range("D49")
$r.Value = "'106.22"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +3.73%  "
$r = $ws.Range("D50")
$r.Value = "'75.34"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$r = $ws.Range("D51")
$r.Value = "'2.778.82"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
